# Atualizacao dos dados bibi add e bdxp
# Updates faturamento_diario data:
#  - Corrects a few values for days in 05/2025 (rows 16, 18, 19, 20)
#  - Inserts a new day (29) for 05/2025 right after day 28, shifting the
#    remaining rows (Apr/Mar/Feb 2025) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct existing values for May/2025 (sheet rows 16, 18, 19, 20)
$ws.Cells.Item(16, 2).Value = 27382.2
$ws.Cells.Item(18, 2).Value = 27829.18
$ws.Cells.Item(19, 2).Value = 15948.71
$ws.Cells.Item(20, 2).Value = 19305.51

# Insert a new row at position 21 (shifts rows 21:79 down to 22:80)
$ws.Rows("21:21").Insert(-4121)

# Fill in the newly inserted row with data for day 29 of May/2025
$ws.Cells.Item(21, 1).Value = 29
$ws.Cells.Item(21, 2).Value = 10955.39
$ws.Cells.Item(21, 3).Value = 5
$ws.Cells.Item(21, 4).Value = 2025
$ws.Cells.Item(21, 5).Value = "05/2025"
